# Updates the "Price" (D) and "Volume(1h)" (E) columns of the cryptos list
# with refreshed values, mirroring a scheduled GitHub Actions data refresh.
# Cells in column D whose new text would otherwise be auto-parsed by Excel
# as a numeric value are prefixed with a leading apostrophe so the cell
# keeps storing the exact literal text (e.g. "141.00", not the number 141),
# matching the original plain-text representation of these cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").Value = "57.912.06"
$ws.Range("E2").Value = "  -2.79%  "
$ws.Range("D3").Value = "2.286.16"
$ws.Range("E3").Value = "  -2.57%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'530.01"
$ws.Range("E5").Value = "  -5.06%  "
$ws.Range("D6").Value = "'131.26"
$ws.Range("E6").Value = "  -0.53%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").Value = "'0.585"
$ws.Range("E8").Value = "  +1.24%  "
$ws.Range("D9").Value = "2.284.63"
$ws.Range("E9").Value = "  -2.59%  "
$ws.Range("D10").Value = "'0.0991"
$ws.Range("E10").Value = "  -4.52%  "
$ws.Range("D11").Value = "'5.46"
$ws.Range("E11").Value = "  -2.29%  "
$ws.Range("E12").Value = "  -0.21%  "
$ws.Range("E13").Value = "  -3.20%  "
$ws.Range("E14").Value = "  -2.15%  "
$ws.Range("D15").Value = "2.692.50"
$ws.Range("E15").Value = "  -2.52%  "
$ws.Range("D16").Value = "57.852.80"
$ws.Range("E16").Value = "  -2.86%  "
$ws.Range("E17").Value = "  -3.52%  "
$ws.Range("D18").Value = "2.284.65"
$ws.Range("E18").Value = "  -2.55%  "
$ws.Range("D19").Value = "'10.49"
$ws.Range("E19").Value = "  -4.31%  "
$ws.Range("D20").Value = "'4.18"
$ws.Range("E20").Value = "  -5.75%  "
$ws.Range("D21").Value = "'310.58"
$ws.Range("E21").Value = "  -2.50%  "
$ws.Range("E22").Value = "  -3.44%  "
$ws.Range("E23").Value = "  -0.06%  "
$ws.Range("D24").Value = "'62.31"
$ws.Range("E24").Value = "  -2.61%  "
$ws.Range("E25").Value = "  -1.93%  "
$ws.Range("E26").Value = "  +0.07%  "
$ws.Range("D27").Value = "'7.96"
$ws.Range("E27").Value = "  -4.33%  "
$ws.Range("E28").Value = "  -6.53%  "
$ws.Range("D29").Value = "'170.54"
$ws.Range("E29").Value = "  -0.42%  "
$ws.Range("E30").Value = "  -5.70%  "
$ws.Range("E31").Value = "  -4.13%  "
$ws.Range("D32").Value = "'5.72"
$ws.Range("E32").Value = "  -3.86%  "
$ws.Range("E33").Value = "  -4.88%  "
$ws.Range("E34").Value = "  -5.25%  "
$ws.Range("E35").Value = "  +0.03%  "
$ws.Range("D36").Value = "'17.72"
$ws.Range("E36").Value = "  -1.55%  "
$ws.Range("E37").Value = "  -0.08%  "
$ws.Range("E38").Value = "  -5.80%  "
$ws.Range("D39").Value = "'3.88"
$ws.Range("E39").Value = "  -4.48%  "
$ws.Range("E40").Value = "  -0.15%  "
$ws.Range("E41").Value = "  -5.23%  "
$ws.Range("D42").Value = "'141.00"
$ws.Range("E42").Value = "  -2.05%  "
$ws.Range("D43").Value = "'285.64"
$ws.Range("E43").Value = "  -8.77%  "
$ws.Range("E44").Value = "  -2.04%  "
$ws.Range("D45").Value = "'0.0946"
$ws.Range("E45").Value = "  -1.14%  "
$ws.Range("E46").Value = "  -2.18%  "
$ws.Range("E47").Value = "  -1.96%  "
$ws.Range("D48").Value = "'17.90"
$ws.Range("E48").Value = "  -5.01%  "
$ws.Range("E49").Value = "  -3.55%  "
$ws.Range("D50").Value = "'10.92"
$ws.Range("E50").Value = "  -1.34%  "
